$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 12-18 (Word table row numbers) represent the schedule entries for
# 01-09-2025 through 07-09-2025. Each of these turns into an overnight
# shift: start time becomes 23:00, end date becomes the next calendar
# day, and end time becomes 07:00. The start date itself is unchanged.
$endDates = @{
    12 = "02-09-2025"
    13 = "03-09-2025"
    14 = "04-09-2025"
    15 = "05-09-2025"
    16 = "06-09-2025"
    17 = "07-09-2025"
    18 = "08-09-2025"
}

foreach ($r in 12..18) {
    $row = $t.Rows.Item($r)
    $row.Cells.Item(4).Range.Text = "23:00"
    $row.Cells.Item(5).Range.Text = $endDates[$r]
    $row.Cells.Item(6).Range.Text = "07:00"
}

# Rows 19-26 (Word table row numbers) represent the schedule entries for
# 08-09-2025 through 18-09-2025. These entries are removed entirely: the
# route, both dates, both times, "V.P." and "36" fields are all cleared.
foreach ($r in 19..26) {
    $row = $t.Rows.Item($r)
    $row.Cells.Item(1).Range.Text = ""
    $row.Cells.Item(3).Range.Text = ""
    $row.Cells.Item(4).Range.Text = ""
    $row.Cells.Item(5).Range.Text = ""
    $row.Cells.Item(6).Range.Text = ""
    $row.Cells.Item(10).Range.Text = ""
    $row.Cells.Item(11).Range.Text = ""
}

Write-Host "edits applied"
